# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# on the cryptos worksheet to reflect the refreshed GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    # Force the cell to store $Text verbatim as text, even when it looks numeric,
    # then restore the default (Normal) style so no visible formatting changes.
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = '35.280.29'
$ws.Range("E2").Value = '  +1.05%  '

# Row 3
$ws.Range("D3").Value = '1.865.36'
$ws.Range("E3").Value = '  +1.35%  '

# Row 4
$ws.Range("E4").Value = '  +0.63%  '

# Row 5
Set-TextValue $ws.Range("D5") '239.81'
$ws.Range("E5").Value = '  +3.52%  '

# Row 6
$ws.Range("E6").Value = '  +0.67%  '

# Row 7
$ws.Range("E7").Value = '  +0.60%  '

# Row 8
Set-TextValue $ws.Range("D8") '42.61'
$ws.Range("E8").Value = '  +6.82%  '

# Row 9
$ws.Range("E9").Value = '  +0.81%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.0695'
$ws.Range("E10").Value = '  +1.22%  '

# Row 11
$ws.Range("E11").Value = '  +0.45%  '

# Row 12
$ws.Range("D12").Value = '2.135.08'
$ws.Range("E12").Value = '  +1.32%  '

# Row 13
Set-TextValue $ws.Range("D13") '11.51'
$ws.Range("E13").Value = '  +0.62%  '

# Row 14
$ws.Range("D14").Value = '1.863.94'
$ws.Range("E14").Value = '  +0.96%  '

# Row 15
Set-TextValue $ws.Range("D15") '0.679'
$ws.Range("E15").Value = '  +0.97%  '

# Row 16
Set-TextValue $ws.Range("D16") '4.73'
$ws.Range("E16").Value = '  +1.78%  '

# Row 17
$ws.Range("D17").Value = '35.284.34'
$ws.Range("E17").Value = '  +1.01%  '

# Row 18
Set-TextValue $ws.Range("D18") '70.09'
$ws.Range("E18").Value = '  +0.28%  '

# Row 19
$ws.Range("D19").Value = '0.0₃0797'
$ws.Range("E19").Value = '  +1.18%  '

# Row 20
Set-TextValue $ws.Range("D20") '241.28'
$ws.Range("E20").Value = '  +0.30%  '

# Row 21
Set-TextValue $ws.Range("D21") '12.26'
$ws.Range("E21").Value = '  +0.79%  '

# Row 22
Set-TextValue $ws.Range("D22") '4.75'
$ws.Range("E22").Value = '  +1.26%  '

# Row 23
$ws.Range("E23").Value = '  +0.52%  '

# Row 24
$ws.Range("E24").Value = '  -0.82%  '

# Row 25
Set-TextValue $ws.Range("D25") '169.62'
$ws.Range("E25").Value = '  -0.96%  '

# Row 26
$ws.Range("E26").Value = '  +24.79%  '

# Row 27
Set-TextValue $ws.Range("D27") '8.13'
$ws.Range("E27").Value = '  +4.25%  '

# Row 28
Set-TextValue $ws.Range("D28") '17.75'

# Row 29
$ws.Range("E29").Value = '  +0.55%  '

# Row 30
$ws.Range("E30").Value = '  +1.86%  '

# Row 31
$ws.Range("E31").Value = '  +0.57%  '

# Row 32
$ws.Range("E32").Value = '  +2.14%  '

# Row 33
Set-TextValue $ws.Range("D33") '1.82'
$ws.Range("E33").Value = '  +27.78%  '

# Row 34
Set-TextValue $ws.Range("D34") '4.04'
$ws.Range("E34").Value = '  +2.11%  '

# Row 35
Set-TextValue $ws.Range("D35") '2.08'
$ws.Range("E35").Value = '  +9.18%  '

# Row 36
Set-TextValue $ws.Range("D36") '0.820'
$ws.Range("E36").Value = '  +17.74%  '

# Row 37
$ws.Range("E37").Value = '  +6.07%  '

# Row 38
$ws.Range("E38").Value = '  +3.72%  '

# Row 39
$ws.Range("E39").Value = '  +4.56%  '

# Row 40
Set-TextValue $ws.Range("D40") '90.50'
$ws.Range("E40").Value = '  -0.79%  '

# Row 41
$ws.Range("D41").Value = '1.345.50'
$ws.Range("E41").Value = '  +0.24%  '

# Row 42
Set-TextValue $ws.Range("D42") '15.26'
$ws.Range("E42").Value = '  +3.03%  '

# Row 43
$ws.Range("E43").Value = '  +15.43%  '

# Row 44
$ws.Range("E44").Value = '  +3.07%  '

# Row 45
$ws.Range("E45").Value = '  +0.61%  '

# Row 46
Set-TextValue $ws.Range("D46") '12.41'
$ws.Range("E46").Value = '  +44.05%  '

# Row 47
$ws.Range("E47").Value = '  -0.91%  '

# Row 48
$ws.Range("E48").Value = '  +4.97%  '

# Row 49
$ws.Range("D49").Value = '2.051.85'
$ws.Range("E49").Value = '  +1.51%  '

# Row 50
$ws.Range("E50").Value = '  +3.29%  '

# Row 51
$ws.Range("E51").Value = '  +1.77%  '

